# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the zh-cn and de-de handback rows.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 05:15:00"
$wsZhCn.Range("H2").Value = "2016-03-24 05:15:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 05:15:07"
$wsDeDe.Range("H2").Value = "2016-03-24 05:15:32"
